# Add a hyperlinked GitHub-repo URL textbox to the "GitHub link" slide
# (slide #33 in presentation order) of the FWD-PPT-video animation deck.
#
# Target shape (from the OOXML diff):
#   <p:sp> "TextBox 2", id=3
#     off  x=1819836  y=1434354   (EMU)
#     ext  cx=6293224 cy=646331   (EMU)
#     noFill, wrap="square" + spAutoFit, paragraph algn="l"
#     single run whose rPr carries an <a:hlinkClick> to the repo URL,
#     and whose text *is* that same URL.
#
# PowerPoint's COM object model expresses shape geometry in points, not
# EMU (1 pt = 12700 EMU), so the EMU numbers from the diff are divided
# down before being handed to AddTextbox / the Left/Top/Width/Height
# properties.

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(33)

$left   = 1819836 / $EMU_PER_PT
$top    = 1434354 / $EMU_PER_PT
$width  = 6293224 / $EMU_PER_PT
$height = 646331 / $EMU_PER_PT

$url = "https://github.com/vaishnavidamodharan800-alt/TNSDC-FWD-VIDEO-ANIMATION-.git"

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 2"

# <a:bodyPr wrap="square" ...><a:spAutoFit/></a:bodyPr>
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1

# <p:spPr> ... <a:noFill/>
$tb.Fill.Visible = 0

$tr = $tb.TextFrame.TextRange
$tr.Text = $url

# <a:pPr algn="l"/>
$tr.ParagraphFormat.Alignment = 1

# <a:rPr ...><a:hlinkClick r:id="rId3"/></a:rPr>
$tr.ActionSettings.Item(1).Hyperlink.Address = $url

Write-Host "Inserted hyperlinked textbox on slide 33: $url"
